# CS342Logbook.xlsx - add another logbook entry mentioning thread_get_priority()
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5, column I ("Function in file" under "Function of main focus in session")
# previously just listed thread_set_priority(); append thread_get_priority() too.
$ws.Range("I5").Value = "thread_set_priority() thread_get_priority() "

# The author's active cell ended up on I6 after making the edit.
$ws.Range("I6").Select()
